# Generate Report for Handoff
# Re-populate the localization-status report: two files moved from
# "handed back" to "ready for handoff" again, under new generated
# filenames. Target-file / handback-file columns (F/G) are no longer
# known at handoff time, so they are cleared; the handback datetime is
# reset to the zero-date sentinel.

$wb = $excel.ActiveWorkbook

$uuid1 = "ce934a49-8368-4a72-bef7-361bc6b7e0e5"
$uuid2 = "ffffc314c950-ec5b-433f-8639-1f57770b872d"
$xlfHash = "cbef67f6a400d9862e042bb21c3fe117ce7570b8"

$status = "Ready for handoff"
$overviewDate = "2016-53-21 00:53:56"
$handbackDt = "0001-01-01 00:00:00"

# per-language "Latest Handoff Datetime" (column E) on the zh-cn / de-de sheets
$handoffDt = @{ "zh-cn" = "2016-03-21 00:53:52"; "de-de" = "2016-03-21 00:53:56" }

function Set-HyperlinkDisplay($ws, $addr, $text) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.TextToDisplay = $text
            break
        }
    }
}

function Remove-HyperlinkAt($ws, $addr) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.Delete()
            break
        }
    }
}

# ---------------------------------------------------------------
# Overview sheet: File Name / zh-cn / de-de / Latest Handoff Date
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$uuid1.md"
Set-HyperlinkDisplay $wsOverview '$A$2' "$uuid1.md"
$wsOverview.Range("B2").Value = $status
$wsOverview.Range("C2").Value = $status
$wsOverview.Range("D2").Value = $overviewDate

$wsOverview.Range("A3").Value = "$uuid2.md"
Set-HyperlinkDisplay $wsOverview '$A$3' "$uuid2.md"
$wsOverview.Range("B3").Value = $status
$wsOverview.Range("C3").Value = $status
$wsOverview.Range("D3").Value = $overviewDate

# ---------------------------------------------------------------
# Per-language detail sheets: zh-cn, de-de
# ---------------------------------------------------------------
$langs = @("zh-cn", "de-de")
foreach ($lang in $langs) {
    $ws = $wb.Worksheets.Item($lang)

    $xlf = "$uuid1.$xlfHash.$lang.xlf"
    $dt = $handoffDt[$lang]

    # Row 2
    $ws.Range("A2").Value = "$uuid1.md"
    Set-HyperlinkDisplay $ws '$A$2' "$uuid1.md"
    $ws.Range("C2").Value = $status
    $ws.Range("D2").Value = $xlf
    Set-HyperlinkDisplay $ws '$D$2' $xlf
    $ws.Range("E2").Value = $dt
    $ws.Range("F2:G2").Clear()
    Remove-HyperlinkAt $ws '$F$2'
    Remove-HyperlinkAt $ws '$G$2'
    $ws.Range("H2").Value = $handbackDt

    # Row 3
    $ws.Range("A3").Value = "$uuid2.md"
    Set-HyperlinkDisplay $ws '$A$3' "$uuid2.md"
    $ws.Range("C3").Value = $status
    $ws.Range("D3").Value = $xlf
    Set-HyperlinkDisplay $ws '$D$3' $xlf
    $ws.Range("E3").Value = $dt
    $ws.Range("F3:G3").Clear()
    Remove-HyperlinkAt $ws '$F$3'
    Remove-HyperlinkAt $ws '$G$3'
    $ws.Range("H3").Value = $handbackDt
}
